# Add header row to Sheet1: Name | type of work | work date
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "type of work"
$ws.Range("C1").Value = "work date"

# Leave the selection on C1, matching the authored workbook's saved cursor position
$ws.Range("C1").Select()
